# final_exercise_day1.docx edit script
#
# 1) Insert a new bold+underlined "Memory File Cache" heading paragraph right
#    after the title paragraph ("Final Exercise - day 1:").
# 2) Remove the stray _GoBack bookmark from the "Use async / await" paragraph.
# 3) Append a whole new "Echo at Time" exercise section at the end of the
#    document (before sectPr), re-adding the _GoBack bookmark at its new
#    location, followed by two empty trailing paragraphs.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Append-Xml([string]$bodyXml) {
    $r = $d.Range($d.Content.End, $d.Content.End)
    $r.InsertXML((New-PkgXml $bodyXml))
}

# ---------------------------------------------------------------------------
# 1) Insert "Memory File Cache" heading paragraph right after the title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item(2)
$headingRange = $headingPara.Range
$headingBody = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Memory File Cache</w:t></w:r></w:p>'
$headingRange.InsertXML((New-PkgXml $headingBody))

# ---------------------------------------------------------------------------
# 2) Remove the _GoBack bookmark that currently sits inside the
#    "Use async / await ..." paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3) Append the new "Echo at Time" section at the end of the document.
# ---------------------------------------------------------------------------

# 3a) Everything from the blank spacer paragraph through the
#     "npm-schedule module." paragraph.
$quoteOpen = [string][char]0x2018
$quoteClose = [string][char]0x2019

$part1 = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p>' +
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Echo at Time</w:t></w:r></w:p>' +
'<w:p><w:r><w:t>Write a server that will print to the console a message at a given time.</w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">The server should have one </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>model ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> called message, with 3 fields:</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Time</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Message</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Status</w:t></w:r></w:p>' +
'<w:p><w:r><w:t>Have a controller that will</w:t></w:r><w:r><w:t xml:space="preserve"> have only he following methods</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Post new message</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Get all messages</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Get message by id </w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">The scheduling </w:t></w:r><w:r><w:t xml:space="preserve">will be done with a timer of your choice. </w:t></w:r></w:p>' +
'<w:p><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">To implement the scheduler use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-schedule module.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

Append-Xml $part1

# 3b) The node-schedule hyperlink paragraph (built via Hyperlinks.Add so the
#     Hyperlink character style survives).
$linkUrl = "https://www.npmjs.com/package/node-schedule"
Append-Xml ('<w:p><w:r><w:t>' + $linkUrl + '</w:t></w:r></w:p>')
$linkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$linkRange = $linkPara.Range
$linkRange.MoveEnd(1, -1)
$d.Hyperlinks.Add($linkRange, $linkUrl) | Out-Null

# 3c) The remaining paragraphs, ending with the _GoBack bookmark re-added,
#     plus the two trailing empty paragraphs.
$part3 = '<w:p><w:r><w:t>Once print is scheduled, change the message status to ' + $quoteOpen + 'Done' + $quoteClose + '.</w:t></w:r></w:p>' +
'<w:p><w:r><w:t xml:space="preserve">Write all your code with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>async</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-await and promises.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
'<w:p/><w:p/>'

Append-Xml $part3

Write-Host "Final paragraph count:" $d.Paragraphs.Count
